$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.378.05'
$ws.Range("D3").Value = '1.882.39'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D5").Value = '''0.7123'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '''243.07'
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("D9").Value = '''0.3158'
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("D10").Value = '''25.10'
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("D11").Value = '''0.08337'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("D12").Value = '1.900.25'
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '''5.271'
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("E14").Value = '  +4.10%  '
$ws.Range("D15").Value = '''0.7186'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '''6.370'
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").Value = '''0.000008668'
$ws.Range("E17").Value = '  +5.31%  '
$ws.Range("D18").Value = '29.406.03'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '''242.95'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '2.160.15'
$ws.Range("E20").Value = '  +1.98%  '
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '''7.845'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("D25").Value = '''0.1576'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").Value = '''9.095'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '''163.40'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '''18.62'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").Value = '''1.510'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '''4.446'
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").Value = '''4.353'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("E32").Value = '  -6.25%  '
$ws.Range("D33").Value = '''0.05393'
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("D34").Value = '''1.949'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").Value = '''0.7750'
$ws.Range("E35").Value = '  +4.02%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = '''2.688'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '1.274.01'
$ws.Range("E39").Value = '  +3.68%  '
$ws.Range("D40").Value = '''2.745'
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("D41").Value = '''6.526'
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Value = '''0.9210'
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("D43").Value = '''113.26'
$ws.Range("E43").Value = '  +2.76%  '
$ws.Range("D44").Value = '''74.53'
$ws.Range("E44").Value = '  +2.41%  '
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("E46").Value = '  +4.92%  '
$ws.Range("D47").Value = '2.045.35'
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("D48").Value = '''1.816'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = '''0.5226'
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").Value = '''9.587'
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").Value = '''0.4383'
$ws.Range("E51").Value = '  +1.22%  '
